# Fix the double space typo in the "PART NO / PART SPEC" header (cell B1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "PART NO / PART SPEC"

# Update the active selection, matching the author's last cursor position on save.
[void]$ws.Range("O12").Select()
